$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 6 trailing rows (36-41) so the sheet shrinks from 41 rows to 35 rows
$ws.Range("A36:C41").EntireRow.Delete()

# Row 2
$ws.Range("A2").Value = "start"
$ws.Range("B2").Clear()
$ws.Range("C2").Value = "AvvioPerizia"

# Row 3
$ws.Range("A3").Value = "start"
$ws.Range("B3").Clear()
$ws.Range("C3").Value = "AvvioPeriziaPostDesk"

# Row 4
$ws.Range("A4").Value = "AvvioPerizia"
$ws.Range("B4").Value = "PeriziaAvviata"
$ws.Range("B4").Borders.LineStyle = 1
$ws.Range("C4").Value = "EsecuzionePerizia"

# Row 5
$ws.Range("A5").Value = "AvvioPeriziaPostDesk"
$ws.Range("B5").Value = "PeriziaAvviata"
$ws.Range("B5").Borders.LineStyle = 1
$ws.Range("C5").Value = "EsecuzionePerizia"

# Row 6
$ws.Range("A6").Value = "EsecuzionePerizia"
$ws.Range("B6").Value = "RevisioneAvviata"
$ws.Range("B6").Borders.LineStyle = 1
$ws.Range("C6").Value = "PeriziaInRevisione"

# Row 7
$ws.Range("A7").Value = "EsecuzionePerizia"
$ws.Range("B7").Value = "AuthorityNonNecessaria"
$ws.Range("B7").Borders.LineStyle = 1
$ws.Range("C7").Value = "ComunicazioneEsito"

# Row 8
$ws.Range("A8").Value = "EsecuzionePerizia"
$ws.Range("B8").Value = "CollabAvviata"
$ws.Range("B8").Borders.LineStyle = 1
$ws.Range("C8").Value = "CollaborazionePerizia"

# Row 9
$ws.Range("A9").Value = "EsecuzionePerizia"
$ws.Range("B9").Value = "RevisioneAvviata"
$ws.Range("B9").Borders.LineStyle = 1
$ws.Range("C9").Value = "RevisionePerizia"

# Row 10
$ws.Range("A10").Value = "EsecuzionePerizia"
$ws.Range("B10").Value = "AuthorityNecessaria"
$ws.Range("B10").Borders.LineStyle = 1
$ws.Range("C10").Value = "VerificaAuthority"

# Row 11
$ws.Range("A11").Value = "PeriziaInRevisione"
$ws.Range("B11").Value = "RevisioniCompletate"
$ws.Range("B11").Borders.LineStyle = 1
$ws.Range("C11").Value = "EsecuzionePerizia"

# Row 12
$ws.Range("A12").Value = "PeriziaInRevisione"
$ws.Range("B12").Value = "RevisioniNonCompletate"
$ws.Range("B12").Borders.LineStyle = 1
$ws.Range("C12").Value = "PeriziaInRevisione"

# Row 13
$ws.Range("A13").Value = "ComunicazioneEsito"
$ws.Range("B13").Value = "EsitoScrittoEmail"
$ws.Range("B13").Borders.LineStyle = 1
$ws.Range("C13").Value = "InvioEmailEsito"

# Row 14
$ws.Range("A14").Value = "ComunicazioneEsito"
$ws.Range("B14").Value = "EsitoScrittoAtto"
$ws.Range("B14").Borders.LineStyle = 1
$ws.Range("C14").Value = "InvioAtto"

# Row 15
$ws.Range("A15").Value = "ComunicazioneEsito"
$ws.Range("B15").Value = "AccordoNonRaggiunto"
$ws.Range("B15").Borders.LineStyle = 1
$ws.Range("C15").Value = "ConciliazionePerizia"

# Row 16
$ws.Range("A16").Value = "ComunicazioneEsito"
$ws.Range("B16").Value = "AccordoRaggiunto-o-NonRichiesto"
$ws.Range("B16").Borders.LineStyle = 1
$ws.Range("C16").Value = "VerificaChiusura"

# Row 17
$ws.Range("A17").Value = "InvioEmailEsito"
$ws.Range("B17").Value = "EmailEsitoInviata"
$ws.Range("B17").Borders.LineStyle = 1
$ws.Range("C17").Value = "AttesaRicezioneAtto"

# Row 18
$ws.Range("A18").Value = "InvioAtto"
$ws.Range("B18").Value = "AttoInviato"
$ws.Range("B18").Borders.LineStyle = 1
$ws.Range("C18").Value = "AttesaRicezioneAtto"

# Row 19
$ws.Range("A19").Value = "AuthorityNonApprovata"
$ws.Range("B19").Value = "PeriziaNonApprovata"
$ws.Range("B19").Borders.LineStyle = 1
$ws.Range("C19").Value = "EsecuzionePerizia"

# Row 20
$ws.Range("A20").Value = "CollaborazionePerizia"
$ws.Range("B20").Value = "CollabRifiutata"
$ws.Range("B20").Borders.LineStyle = 1
$ws.Range("C20").Value = "CollaborazioneRifiutata"

# Row 21
$ws.Range("A21").Value = "CollaborazionePerizia"
$ws.Range("B21").Value = "CollabCancellata"
$ws.Range("B21").Borders.LineStyle = 1
$ws.Range("C21").Value = "CollaborazioneCancellata"

# Row 22
$ws.Range("A22").Value = "CollaborazionePerizia"
$ws.Range("B22").Value = "CollabCompletata"
$ws.Range("B22").Borders.LineStyle = 1
$ws.Range("C22").Value = "CollaborazioneCompletata"

# Row 23
$ws.Range("A23").Value = "CollaborazioneRifiutata"
$ws.Range("B23").Clear()
$ws.Range("C23").Value = "EsecuzionePerizia"

# Row 24
$ws.Range("A24").Value = "CollaborazioneCancellata"
$ws.Range("B24").Clear()
$ws.Range("C24").Value = "EsecuzionePerizia"

# Row 25
$ws.Range("A25").Value = "CollaborazioneCompletata"
$ws.Range("B25").Clear()
$ws.Range("C25").Value = "EsecuzionePerizia"

# Row 26
$ws.Range("A26").Value = "RevisionePerizia"
$ws.Range("B26").Value = "RevisioneRifiutata"
$ws.Range("B26").Borders.LineStyle = 1
$ws.Range("C26").Value = "RevisioneRifiutata"

# Row 27
$ws.Range("A27").Value = "RevisionePerizia"
$ws.Range("B27").Value = "RevisioneCompletata"
$ws.Range("B27").Borders.LineStyle = 1
$ws.Range("C27").Value = "RevisioneCompletata"

# Row 28
$ws.Range("A28").Value = "RevisioneRifiutata"
$ws.Range("B28").Clear()
$ws.Range("C28").Value = "PeriziaInRevisione"

# Row 29
$ws.Range("A29").Value = "RevisioneCompletata"
$ws.Range("B29").Clear()
$ws.Range("C29").Value = "PeriziaInRevisione"

# Row 30
$ws.Range("A30").Value = "VerificaAuthority"
$ws.Range("B30").Value = "PeriziaApprovata"
$ws.Range("B30").Borders.LineStyle = 1
$ws.Range("C30").Value = "ComunicazioneEsito"

# Row 31
$ws.Range("A31").Value = "VerificaAuthority"
$ws.Range("B31").Clear()
$ws.Range("C31").Value = "AuthorityNonApprovata"

# Row 32
$ws.Range("A32").Value = "ConciliazionePerizia"
$ws.Range("B32").Value = "PeriziaConcilazioneAvviata"
$ws.Range("B32").Borders.LineStyle = 1
$ws.Range("C32").Value = "EsecuzionePeriziaConciliazione"

# Row 33
$ws.Range("A33").Value = "ConciliazionePerizia"
$ws.Range("B33").Value = "ChiusuraConAccordo"
$ws.Range("B33").Borders.LineStyle = 1
$ws.Range("C33").Value = "VerificaChiusura"

# Row 34
$ws.Range("A34").Value = "ConciliazionePerizia"
$ws.Range("B34").Value = "ChiusuraSenzaAccordo"
$ws.Range("B34").Borders.LineStyle = 1
$ws.Range("C34").Value = "VerificaChiusura"

# Row 35
$ws.Range("A35").Value = "VerificaChiusura"
$ws.Range("B35").Value = "PeriziaIncompletaAvviata"
$ws.Range("B35").Borders.LineStyle = 1
$ws.Range("C35").Value = "EsecuzionePerizia"
